$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 2-4 (columns B and D)
$ws.Range("B2").Value = 167.5429929727877
$ws.Range("D2").Value = 0.176811018107935

$ws.Range("B3").Value = 31331.46756994809
$ws.Range("D3").Value = 0.01017445002450822

$ws.Range("B4").Value = 23995.74187856295
$ws.Range("D4").Value = 0.1047600572820151

# Add new rows 5 and 6 - copy formatting of the existing label cell (A4) first,
# then overwrite with the new label text
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "googlesearch_media_cost"
$ws.Range("B5").Value = 44668.42512772494
$ws.Range("C5").Value = 1494447.761988
$ws.Range("D5").Value = 0.02988958614940441

$ws.Range("A4").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "dv360_media_cost"
$ws.Range("B6").Value = 477.8417833700823
$ws.Range("C6").Value = 271129.18
$ws.Range("D6").Value = 0.001762413707628527
